$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell (H1) onto the new header cells I1:J1
# so they reuse the same bold/border/centered style, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add new data values for columns I and J, rows 2-4 (plain, unstyled numeric cells)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8
